$d = $word.ActiveDocument

$r0 = $d.Content
$found0 = $r0.Find.Execute("Appendix 17: SWIFT Referrals Interview: Information Sheet and Consent Form", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found0) { throw "Not found: Appendix 17: SWIFT Referrals Interview: Information Sheet and Consent Form" }
$r0.Text = "Bylaag 17: SWIFT Verwysingsonderhoud: Inligtingsblad en Toestemmingsvorm"

$r1 = $d.Content
$found1 = $r1.Find.Execute(" What will my interview look like and what is expected of me?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found1) { throw "Not found:  What will my interview look like and what is expected of me?" }
$r1.Text = " Hoe sal my onderhoud lyk en wat word van my verwag?"

$r2 = $d.Content
$found2 = $r2.Find.Execute("Why have I been invited to the interview?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found2) { throw "Not found: Why have I been invited to the interview?" }
$r2.Text = "Waarom is ek na die onderhoud genooi?"

$r3 = $d.Content
$found3 = $r3.Find.Execute("Do I have to agree to be interviewed?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found3) { throw "Not found: Do I have to agree to be interviewed?" }
$r3.Text = "Moet ek instem om ondervra te word?"

$r4 = $d.Content
$found4 = $r4.Find.Execute("What happens with my information?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found4) { throw "Not found: What happens with my information?" }
$r4.Text = "Wat gebeur met my inligting?"

$r5 = $d.Content
$found5 = $r5.Find.Execute("To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found5) { throw "Not found: To protect your personal information (including your real name, contact details, and any other information that can identify you), we will give you a participant number, and you can choose a name you want us to call you during the interview." }
$r5.Text = "Om jou persoonlike inligting (insluitend jou regte naam, kontakbesonderhede, en enige ander inligting wat jou kan identifiseer) te beskerm, sal ons vir jou 'n deelnemernommer gee, en jy kan 'n naam kies waarmee ons jou tydens die onderhoud kan aanspreek."

$r6 = $d.Content
$found6 = $r6.Find.Execute("Do I get anything for being interviewed? ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found6) { throw "Not found: Do I get anything for being interviewed? " }
$r6.Text = "Kry ek enige iets vir deelname aan die onderhoud? "

$r7 = $d.Content
$found7 = $r7.Find.Execute("As a thank you for taking part in the discussion, we will give you a R30 airtime voucher/data bundle. ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found7) { throw "Not found: As a thank you for taking part in the discussion, we will give you a R30 airtime voucher/data bundle. " }
$r7.Text = "As 'n bedanking vir jou deelname aan die gesprek, sal ons vir jou 'n R30 lugtydkoepon/data-bundel gee. "

$r8 = $d.Content
$found8 = $r8.Find.Execute("What happens to my information if I agree to be interviewed?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found8) { throw "Not found: What happens to my information if I agree to be interviewed?" }
$r8.Text = "Wat gebeur met my inligting as ek instem om ondervra te word?"

$r9 = $d.Content
$found9 = $r9.Find.Execute("Who are some of the study team members?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found9) { throw "Not found: Who are some of the study team members?" }
$r9.Text = "Wie is sommige van die spanlede van die studie?"

$r10 = $d.Content
$found10 = $r10.Find.Execute("The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found10) { throw "Not found: The principal investigators of this study are Prof Cathy Ward and Cindee Bruyns and the Co-investigator is Carly Katzef all from the University of Cape Town." }
$r10.Text = "Die hoofondersoeker van hierdie studie is Prof. Cathy Ward en Cindee Bruyns, en die Mede-ondersoeker is Carly Katzef, almal van die Universiteit van Kaapstad."

$r11 = $d.Content
$found11 = $r11.Find.Execute("Are there any risks in being interviewed?   ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found11) { throw "Not found: Are there any risks in being interviewed?   " }
$r11.Text = "Is daar enige risiko's verbonde aan die onderhoud?   "

$r12 = $d.Content
$found12 = $r12.Find.Execute("Who pays for the study?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found12) { throw "Not found: Who pays for the study?" }
$r12.Text = "Wie betaal vir die studie?"

$r13 = $d.Content
$found13 = $r13.Find.Execute("This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found13) { throw "Not found: This study is part of the Global Parenting Initiative, funded by the LEGO Foundation, Oak Foundation, the World Childhood Foundation, The Human Safety Net, and the UK Research and Innovation Global Challenges Research Fund. " }
$r13.Text = "Hierdie studie is deel van die Global Parenting Initiative, gefinansier deur die LEGO Foundation, Oak Foundation, die World Childhood Foundation, The Human Safety Net, en die UK Research and Innovation Global Challenges Research Fund. "

$r14 = $d.Content
$found14 = $r14.Find.Execute("Data protection", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found14) { throw "Not found: Data protection" }
$r14.Text = "Databeskerming"

$r15 = $d.Content
$found15 = $r15.Find.Execute("Who has approved this study?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found15) { throw "Not found: Who has approved this study?" }
$r15.Text = "Wie het hierdie studie goedgekeur?"

$r16 = $d.Content
$found16 = $r16.Find.Execute("Who do I contact if I have questions or concerns?", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found16) { throw "Not found: Who do I contact if I have questions or concerns?" }
$r16.Text = "Wie kan ek kontak as ek vrae of bekommernisse het?"

$r17 = $d.Content
$found17 = $r17.Find.Execute("If you have any questions or concerns about your rights as a study participant, you can contact the study team at swift@globalparenting.org or on WhatsApp at +27 XX XXX XXXX (messages only).", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found17) { throw "Not found: If you have any questions or concerns about your rights as a study participant, you can contact the study team at swift@globalparenting.org or on WhatsApp at +27 XX XXX XXXX (messages only)." }
$r17.Text = "As jy enige vrae of bekommernisse het oor jou regte as 'n studie-deelnemer, kan jy die studiespan kontak by swift@globalparenting.org of via WhatsApp by +27 XX XXX XXXX (net boodskappe)."

$r18 = $d.Content
$found18 = $r18.Find.Execute("If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found18) { throw "Not found: If you have more questions or concerns about your rights, you can contact one of the ethics committees listed: " }
$r18.Text = "As jy meer vrae of bekommernisse het oor jou regte, kan jy een van die etiekkomitees hieronder kontak: "

$r19 = $d.Content
$found19 = $r19.Find.Execute("Name", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found19) { throw "Not found: Name" }
$r19.Text = "Naam"

$r20 = $d.Content
$found20 = $r20.Find.Execute("Telephone", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found20) { throw "Not found: Telephone" }
$r20.Text = "Telefoon"

$r21 = $d.Content
$found21 = $r21.Find.Execute("Email", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found21) { throw "Not found: Email" }
$r21.Text = "E-pos"

$r22 = $d.Content
$found22 = $r22.Find.Execute("University of Cape Town Centre for Social Science Research ", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found22) { throw "Not found: University of Cape Town Centre for Social Science Research " }
$r22.Text = "Universiteit van Kaapstad Sentrum vir Sosiale Wetenskap Navorsing "

$r23 = $d.Content
$found23 = $r23.Find.Execute("Human Research Ethics Committee", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found23) { throw "Not found: Human Research Ethics Committee" }
$r23.Text = "Etiekkomitee vir Menslike Navorsing"

$r24 = $d.Content
$found24 = $r24.Find.Execute("Informed Telephonic consent to take part in the study.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found24) { throw "Not found: Informed Telephonic consent to take part in the study." }
$r24.Text = "Ingeligte Telefoniese toestemming om aan die studie deel te neem."

$r25 = $d.Content
$found25 = $r25.Find.Execute("Someone from the research team has gone over all the information above and I know what I need to do.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found25) { throw "Not found: Someone from the research team has gone over all the information above and I know what I need to do." }
$r25.Text = "Iemand van die navorsingspan het al die inligting hierbo deurgegaan en ek weet wat ek moet doen."

$r26 = $d.Content
$found26 = $r26.Find.Execute("I know who can see my information after the interview, how it will be kept safe, and what happens to it after the study.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found26) { throw "Not found: I know who can see my information after the interview, how it will be kept safe, and what happens to it after the study." }
$r26.Text = "Ek weet wie my inligting na die onderhoud kan sien, hoe dit veilig gehou sal word, en wat daarmee sal gebeur na die studie."

$r27 = $d.Content
$found27 = $r27.Find.Execute("I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found27) { throw "Not found: I know I can request access to my data, correct any mistakes, ask to delete it, or for it to be transferred somewhere else." }
$r27.Text = "Ek weet ek kan toegang tot my data versoek, enige foute regstel, vra dat dit verwyder word, of vir dit om na 'n ander plek oorgedra te word."

$r28 = $d.Content
$found28 = $r28.Find.Execute("I know that I won’t be named in any papers or reports from this study.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found28) { throw "Not found: I know that I won’t be named in any papers or reports from this study." }
$r28.Text = "Ek weet dat ek nie in enige artikels of verslae van hierdie studie genoem sal word nie."

$r29 = $d.Content
$found29 = $r29.Find.Execute("I know who to tell if I have a problem with the study.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found29) { throw "Not found: I know who to tell if I have a problem with the study." }
$r29.Text = "Ek weet wie ek moet kontak as ek 'n probleem met die studie het."

$r30 = $d.Content
$found30 = $r30.Find.Execute("I can be contacted again if more information is needed from me.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found30) { throw "Not found: I can be contacted again if more information is needed from me." }
$r30.Text = "Ek kan weer gekontak word as meer inligting van my nodig is."

$r31 = $d.Content
$found31 = $r31.Find.Execute("I understand the team will keep my contact information safe so they can tell me about the results of the study.", $true, $false, $false, $false, $false, $true, 1, $false)
if (-not $found31) { throw "Not found: I understand the team will keep my contact information safe so they can tell me about the results of the study." }
$r31.Text = "Ek verstaan dat die span my kontakbesonderhede veilig sal hou sodat hulle my kan inlig oor die resultate van die studie."

